$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 44, shifting rows 44:48 down to 45:49.
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the new weekly price record.
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").Value = 44476
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 100114001
$ws.Range("G44").Value = "Papa"
$ws.Range("H44").Value = "Rosara"
$ws.Range("I44").Value = "1a (guarda)"
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 11000
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = 11500
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Región del Maule"
$ws.Range("P44").Value = 460
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
